# FINFLUX-3612 Carpetas (Cartias) specific scenarios
# Updates the Summary / Original Schedule / Repayment schedule figures that
# shift once the compounding penalty charge is recalculated for the 1st
# repayment scenario, and moves the active selection on to the
# "Repayment schedule" tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A3").Value = 260.2
$wsSummary.Range("E3").Value = 209.19
$wsSummary.Range("A5").Value = 0
$wsSummary.Range("E5").Value = 0

# ---------------------------------------------------------------------
# Original Schedule sheet
# ---------------------------------------------------------------------
$wsOriginal = $wb.Worksheets.Item("Original Schedule")
$wsOriginal.Range("F4").Value = 0
$wsOriginal.Range("G4").Value = 887.72

# ---------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------
$wsRepayment = $wb.Worksheets.Item("Repayment schedule")

# New column widths for the "Paid Date" helper columns (C:D)
$wsRepayment.Range("C1:D1").EntireColumn.ColumnWidth = 9.140625

$wsRepayment.Range("J4").Value = 0
$wsRepayment.Range("K4").Value = 887.72
$wsRepayment.Range("Q4").Value = 887.72

$wsRepayment.Range("F5").Value = 844.84
$wsRepayment.Range("G5").Value = 2474.71
$wsRepayment.Range("H5").Value = 42.88

$wsRepayment.Range("F6").Value = 844.4
$wsRepayment.Range("G6").Value = 1630.31
$wsRepayment.Range("H6").Value = 43.32

$wsRepayment.Range("F7").Value = 848.2
$wsRepayment.Range("G7").Value = 782.11
$wsRepayment.Range("H7").Value = 39.52

$wsRepayment.Range("F8").Value = 782.11
$wsRepayment.Range("H8").Value = 39.49
$wsRepayment.Range("K8").Value = 821.6
$wsRepayment.Range("Q8").Value = 821.6

# ---------------------------------------------------------------------
# Selections: move through the sheets in tab order so the last one
# selected ("Repayment schedule") ends up as the active tab, matching
# the new activeTab/tabSelected state, while leaving the requested
# activeCell/sqref behind on each sheet along the way.
# ---------------------------------------------------------------------
$wsSummary.Range("D8").Select()
$wsOriginal.Range("D13").Select()
$wsRepayment.Range("K10").Select()
